$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "This test is to get description components.`n"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "0"

$ws.Range("C2").Value = "<class 'AutomationFramework.page_objects.hardware.hardware.Hardware'>"
$ws.Range("D2").Value = "hw_component_description"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = ""

$ws.Range("F2").NumberFormat = "@"
$ws.Range("F2").Value = ""

$ws.Range("G2").Value = "<get>`n  <filter>`n    <components xmlns=""http://openconfig.net/yang/platform"">`n      <component>`n        <name>Waveserver-Ai</name>`n        <state>`n          <description></description>`n        </state>`n      </component>`n    </components>`n  </filter>`n</get>"

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = ""

$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = ""

$ws.Range("J2").Value = "<?xml version=""1.0"" encoding=""UTF-8""?>`n<rpc-reply message-id=""urn:uuid:00f5907b-1911-4318-bcaf-02aa8e8d4834""`n xmlns:ncx=""http://netconfcentral.org/ns/yuma-ncx""`n ncx:last-modified=""2020-10-07T13:51:28Z"" ncx:etag=""814""`n xmlns=""urn:ietf:params:xml:ns:netconf:base:1.0"">`n <data>`n  <components xmlns=""http://openconfig.net/yang/platform"">`n   <component>`n    <name>Waveserver-Ai</name>`n    <state>`n     <description>Waveserver Ai Chassis 3-slot, 1RU</description>`n    </state>`n   </component>`n  </components>`n </data>`n</rpc-reply>"
